$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$check = [char]0x2713

# --- Row 2: Llanelli Town AFC - The New Saints (result recorded) ---
$ws.Range("A2").Value = "Llanelli Town AFC - The New Saints " + $check + ": 0:4"
$ws.Range("C2").Value = 71
$ws.Range("G2").Value = $check

# --- Row 3 becomes SC Uniao Torreense (moved up, result recorded) ---
$ws.Range("A3").Value = "SC Uni" + [char]0xE3 + "o Torreense " + $check + " - UD Oliveirense: 3:2"
$ws.Range("B3").Value = "SC Uni" + [char]0xE3 + "o Torreense"
$ws.Range("C3").Value = 64
$ws.Range("D3").ClearContents()
$ws.Range("E3").Value = 100
$ws.Range("F3").Value = 1.99
$ws.Range("G3").Value = $check

# --- Row 4 becomes AC Milan (moved down, score updated) ---
$ws.Range("A4").Value = "AC Milan  - Pisa Sporting Club: 2:2"
$ws.Range("B4").Value = "AC Milan"
$ws.Range("C4").Value = 56
$ws.Range("D4").Value = 100
$ws.Range("E4").Value = 92
$ws.Range("F4").Value = 1.33
$ws.Range("G4").ClearContents()

# --- Row 5: new fixture HNK Gorica - HNK Hajduk Split ---
$ws.Range("A5").Value = "HNK Gorica - HNK Hajduk Split " + $check + ": 1:3"
$ws.Range("B5").Value = "HNK Hajduk Split"
$ws.Range("C5").Value = 55
$ws.Range("D5").Value = 79
$ws.Range("F5").Value = 1.91
$ws.Range("G5").Value = $check

# --- Row 6: new fixture Aarhus GF - FC Nordsjaelland ---
$ws.Range("A6").Value = "Aarhus GF " + $check + " - FC Nordsjaelland: 1:0"
$ws.Range("B6").Value = "Aarhus GF"
$ws.Range("C6").Value = 52
$ws.Range("D6").Value = 88
$ws.Range("F6").Value = 1.78
$ws.Range("G6").Value = $check
